$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46016
$ws.Range("B2").Value = 33.18
$ws.Range("C2").Value = 30.16
$ws.Range("D2").Value = 23.39
$ws.Range("E2").Value = 15.32
$ws.Range("F2").Value = 10.62
$ws.Range("G2").Value = 15.97
$ws.Range("H2").Value = 21.84
$ws.Range("I2").Value = 38
$ws.Range("J2").Value = 42.64
$ws.Range("K2").Value = 38.97
$ws.Range("L2").Value = 30.72
$ws.Range("M2").Value = 24.28
$ws.Range("N2").Value = 31.46
$ws.Range("O2").Value = 30.98
$ws.Range("P2").Value = 23.23
$ws.Range("Q2").Value = 35.5
$ws.Range("R2").Value = 63.21
$ws.Range("S2").Value = 89.51000000000001
$ws.Range("T2").Value = 95.77
$ws.Range("U2").Value = 93.27
$ws.Range("V2").Value = 97.03
$ws.Range("W2").Value = 99.17
$ws.Range("X2").Value = 94
$ws.Range("Y2").Value = 90.56999999999999
$ws.Range("Z2").Value = 48.7
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 95.19
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 98.09999999999999
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 94.52
$ws.Range("AG2").Value = "0h-15h"
